# Applies the cryptos-list price/volume/coin updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.863.01"
$ws.Range("E2").Value = "  -0.57%  "

# Row 3
$ws.Range("D3").Value = "'3.537.29"
$ws.Range("E3").Value = "  -2.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.25%  "

# Row 5
$ws.Range("D5").Value = "'197.20"
$ws.Range("E5").Value = "  -2.98%  "

# Row 6
$ws.Range("D6").Value = "'555.72"
$ws.Range("E6").Value = "  -1.70%  "

# Row 7
$ws.Range("D7").Value = "'0.658"
$ws.Range("E7").Value = "  +6.63%  "

# Row 8
$ws.Range("D8").Value = "'3.529.67"
$ws.Range("E8").Value = "  -2.38%  "

# Row 9
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").Value = "'0.663"
$ws.Range("E10").Value = "  -1.99%  "

# Row 11
$ws.Range("D11").Value = "'60.80"
$ws.Range("E11").Value = "  +5.14%  "

# Row 12
$ws.Range("E12").Value = "  -5.60%  "

# Row 13
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  -6.90%  "

# Row 14
$ws.Range("D14").Value = "'9.94"
$ws.Range("E14").Value = "  -0.77%  "

# Row 15
$ws.Range("D15").Value = "'4.108.02"
$ws.Range("E15").Value = "  -2.52%  "

# Row 16
$ws.Range("D16").Value = "'3.544.39"
$ws.Range("E16").Value = "  -2.54%  "

# Row 17
$ws.Range("E17").Value = "  -1.62%  "

# Row 18
$ws.Range("D18").Value = "'67.732.33"
$ws.Range("E18").Value = "  -0.60%  "

# Row 19
$ws.Range("D19").Value = "'18.42"
$ws.Range("E19").Value = "  -0.61%  "

# Row 20
$ws.Range("D20").Value = "'11.92"
$ws.Range("E20").Value = "  -4.19%  "

# Row 21
$ws.Range("D21").Value = "'1.03"
$ws.Range("E21").Value = "  -4.64%  "

# Row 22
$ws.Range("D22").Value = "'400.07"
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("D23").Value = "'4.01"
$ws.Range("E23").Value = "  -3.98%  "

# Row 24
$ws.Range("D24").Value = "'86.98"
$ws.Range("E24").Value = "  +1.70%  "

# Row 25
$ws.Range("D25").Value = "'11.77"
$ws.Range("E25").Value = "  -9.16%  "

# Row 26
$ws.Range("D26").Value = "'12.44"
$ws.Range("E26").Value = "  -0.95%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'3.87"
$ws.Range("E27").Value = "  +0.96%  "

# Row 28
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'2.84"
$ws.Range("E28").Value = "  -3.71%  "

# Row 29
$ws.Range("D29").Value = "'8.93"
$ws.Range("E29").Value = "  -2.24%  "

# Row 30
$ws.Range("D30").Value = "'722.76"
$ws.Range("E30").Value = "  +4.28%  "

# Row 31
$ws.Range("D31").Value = "'31.38"
$ws.Range("E31").Value = "  -1.33%  "

# Row 32
$ws.Range("D32").Value = "'7.07"
$ws.Range("E32").Value = "  -13.34%  "

# Row 33
$ws.Range("D33").Value = "'11.79"
$ws.Range("E33").Value = "  -3.46%  "

# Row 34
$ws.Range("D34").Value = "'64.43"
$ws.Range("E34").Value = "  +0.32%  "

# Row 35
$ws.Range("E35").Value = "  -3.28%  "

# Row 36
$ws.Range("D36").Value = "'38.69"
$ws.Range("E36").Value = "  -8.96%  "

# Row 37
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").Value = "'0.393"
$ws.Range("E38").Value = "  -7.15%  "

# Row 39
$ws.Range("D39").Value = "'0.132"
$ws.Range("E39").Value = "  -4.93%  "

# Row 40
$ws.Range("D40").Value = "'3.01"
$ws.Range("E40").Value = "  -3.50%  "

# Row 41
$ws.Range("D41").Value = "'3.080.43"
$ws.Range("E41").Value = "  -5.07%  "

# Row 42
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.06%  "

# Row 43
$ws.Range("D43").Value = "'0.0₃0687"
$ws.Range("E43").Value = "  -11.49%  "

# Row 44
$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = "  -10.61%  "

# Row 45
$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = "  +1.82%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.136"
$ws.Range("E46").Value = "  +3.53%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0410"
$ws.Range("E47").Value = "  -1.89%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.04"
$ws.Range("E48").Value = "  -1.64%  "

# Row 49
$ws.Range("D49").Value = "'139.47"
$ws.Range("E49").Value = "  -1.36%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'2.57"
$ws.Range("E50").Value = "  -14.51%  "

# Row 51
$ws.Range("D51").Value = "'8.29"
$ws.Range("E51").Value = "  -6.87%  "
